# Daily auto push update: insert two new rows of data right before the
# existing "2026/12/29" block (old row 616), shifting all subsequent rows
# down by 2. The two new rows continue the "2026/01/11" / "2026/01/12"
# date sequence that the sheet was missing.
#
#   new row 616: 2026/01/11  日  23  201
#   new row 617: 2026/01/12  月   1  201

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 616 (pushes old row 616.. down to 618..)
$ws.Rows.Item(616).Resize(2).Insert()

# --- Row 616 ---------------------------------------------------------
$cellA616 = $ws.Cells.Item(616, 1)
# Force text formatting first so Excel doesn't auto-convert the
# "yyyy/mm/dd" looking string into a date serial number.
$cellA616.NumberFormat = "@"
$cellA616.Value = "2026/01/11"
# Restore the default (unstyled) look used by all the other data rows.
$cellA616.Style = "Normal"

$ws.Cells.Item(616, 2).Value = "日"
$ws.Cells.Item(616, 3).Value = 23
$ws.Cells.Item(616, 4).Value = 201

# --- Row 617 ---------------------------------------------------------
$cellA617 = $ws.Cells.Item(617, 1)
$cellA617.NumberFormat = "@"
$cellA617.Value = "2026/01/12"
$cellA617.Style = "Normal"

$ws.Cells.Item(617, 2).Value = "月"
$ws.Cells.Item(617, 3).Value = 1
$ws.Cells.Item(617, 4).Value = 201
